$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 76 values (B..W)
$row76 = @(4902.965208700517,909.77832594394704,-196.35657995417239,-260.95010508196384,-361.32977743852831,3501.8361177763713,-402.82888923917517,-131.85398300954628,-30.399506203604364,-55.317171646238855,-340.32066931272715,-92.37166566066287,-24.32915355370805,-62.111510173449886,668.56844620231789,-94.521377930635964,18.32991108093951,84.367770846312268,315.61835965303032,-61.323159791283821,-12.852829714954414,-1.2848959438638907)

# New row 77 values (B..W)
$row77 = @(1.4757616550993236,2.1772607417544743,3.0434831170432886,1.2415832626413441,1.2457268877900551,0.56126587903722047,1.1743272631947785,0.23054550741241203,-3.4261782347163527,0.66334587487468022,1.0306370471546036,0.76256777522881058,1.2988415213542728,2.4136597060060243,1.2247158941234109,131.36807363693123,0.50972075458921062,1.623167155359619,2.1677719788989656,-0.26018393796375672,6.8111669711356573,-2.5626701301248795)

# Column A label for the two new rows
$ws.Range("A76").Value = "negative"
$ws.Range("A77").Value = "negative"

for ($i = 0; $i -lt $row76.Length; $i++) {
    $col = $i + 2  # column B is index 2
    $ws.Cells.Item(76, $col).Value = $row76[$i]
    $ws.Cells.Item(77, $col).Value = $row77[$i]
}

$ws.Range("A1:W77").Select()
